# Apply the "Кеклеон (фиолетовый)" dialogue-table update:
#  - restyle existing row 14 to the "group divider" look (A14 added, s=6/7)
#  - append 9 new rows (15-23) of English / Russian / "converted" dialogue text
#  - move the view's active cell to D19 (matches the new data range)
#
# New shared-string text is written column-by-column within each dialogue
# group (C's first, then D's, then E's) so the shared-string table ends up
# with the same index ordering (45-71) as the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Row 14: turn it into a "last row of a dialogue group" row - copy the
#    border/font formatting used by rows 3 / 16 / 19 / 21 (A:E, style 6/7)
#    onto row 14. This keeps B14/C14/D14/E14 values untouched and adds the
#    empty, styled A14 cell.
# ---------------------------------------------------------------------
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A14:E14").PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------
# 2) Lay down formatting + numbers for the new rows 15-23 first (this
#    does not touch the shared-strings table).
# ---------------------------------------------------------------------
$ws.Range("B4:E4").Copy() | Out-Null
$ws.Range("B15:E15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B15").Value = 137
$ws.Rows.Item(15).RowHeight = 52.2

$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A16:E16").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B16").Value = 140
$ws.Rows.Item(16).RowHeight = 21.6

$ws.Range("B4:E4").Copy() | Out-Null
$ws.Range("B17:E17").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B17").Value = 79

$ws.Range("B4:E4").Copy() | Out-Null
$ws.Range("B18:E18").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B18").Value = 82
$ws.Rows.Item(18).RowHeight = 31.8

$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A19:E19").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B19").Value = 85
$ws.Rows.Item(19).RowHeight = 31.8

$ws.Range("B4:E4").Copy() | Out-Null
$ws.Range("B20:E20").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B20").Value = 52
$ws.Rows.Item(20).RowHeight = 30.6

$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A21:E21").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B21").Value = 55
$ws.Rows.Item(21).RowHeight = 21.6

$ws.Range("B4:E4").Copy() | Out-Null
$ws.Range("B22:E22").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B22").Value = 25
$ws.Rows.Item(22).RowHeight = 21.6

$ws.Range("B4:E4").Copy() | Out-Null
$ws.Range("B23:E23").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B23").Value = 28
$ws.Rows.Item(23).RowHeight = 31.8

# ---------------------------------------------------------------------
# 3) Now write the dialogue text itself, grouped & ordered exactly as the
#    original authoring tool produced the shared-strings table: within
#    each dialogue group, all English (C) cells, then all Russian (D)
#    cells, then all "converted" (E) cells.
# ---------------------------------------------------------------------

# Group: rows 15-16
$ws.Range("C15").Value = " My brother is waiting for a\nchance to speed over to see Team [CS:X]Charm[CR].\nIt\'s quite annoying!"
$ws.Range("C16").Value = " For when that very chance\narrives...[K] I intend to dash first!"
$ws.Range("D15").Value = " Мой брат очень хочет поскорее\nулизнуть из магазина, чтобы повидать\nКоманду [CS:X]Шарм[CR]. Как это меня раздражает!"
$ws.Range("D16").Value = " Когда представится возможность...[K]\nЯ побегу к ним первым!"
$ws.Range("E15").Value = " Íïê áñàó ïœåîû öïœåó ðïòëïñåå\nôìéèîôóû éè íàãàèéîà, œóïáú ðïâéäàóû\nËïíàîäô [CS:X]Šàñí[CR]. Ëàë üóï íåîÿ ñàèäñàçàåó!"
$ws.Range("E16").Value = " Ëïãäà ðñåäòóàâéóòÿ âïèíïçîïòóû...[K]\nŸ ðïáåãô ë îéí ðåñâúí!"

# Group: rows 17-18-19
$ws.Range("C17").Value = " ...Hmmm…"
$ws.Range("C18").Value = " What could it possibly be? What\nhas made little [CS:N]Azurill[CR] unable to wake?"
$ws.Range("C19").Value = " I do hope the dear child will\nawaken and brighten our day with his cheer."
$ws.Range("D17").Value = " ...Хммм..."
$ws.Range("D18").Value = " Что же это может быть? Почему\nмалютка [CS:N]Азурилл[CR] не может проснуться?"
$ws.Range("D19").Value = " Надеюсь, что малыш проснётся и\nснова обрадует нас своей жизнерадостностью."
$ws.Range("E17").Value = " ...Öííí..."
$ws.Range("E18").Value = " Œóï çå üóï íïçåó áúóû? Ðïœåíô\níàìýóëà [CS:N]Àèôñéìì[CR] îå íïçåó ðñïòîôóûòÿ?"
$ws.Range("E19").Value = " Îàäåýòû, œóï íàìúš ðñïòîæóòÿ é\nòîïâà ïáñàäôåó îàò òâïåê çéèîåñàäïòóîïòóûý."

# Group: rows 20-21
$ws.Range("C20").Value = " Please don\'t succumb to\nanyone bad!"
$ws.Range("C21").Value = " Go forward to victory! Fight!\nFight!"
$ws.Range("D20").Value = " Прошу, не уступите негодяям!"
$ws.Range("D21").Value = " Вперёд, к победе! Сражайтесь!\nБоритесь!"
$ws.Range("E20").Value = " Ðñïšô, îå ôòóôðéóå îåãïäÿÿí!"
$ws.Range("E21").Value = " Âðåñæä, ë ðïáåäå! Òñàçàêóåòû!\nÁïñéóåòû!"

# Group: rows 22-23
$ws.Range("C22").Value = " Ah, the returning heroes\nof Team [team:]!"
$ws.Range("C23").Value = " I thank you sincerely for\nstopping the spreading nightmare!"
$ws.Range("D22").Value = " Ах, наши герои из Команды\n[team:]!"
$ws.Range("D23").Value = " Я от всего сердца благодарю\nвас, что вы остановили надвигающийся\nкошмар!"
$ws.Range("E22").Value = " Àö, îàšé ãåñïé éè Ëïíàîäú\n[team:]!"
$ws.Range("E23").Value = " Ÿ ïó âòåãï òåñäøà áìàãïäàñý\nâàò, œóï âú ïòóàîïâéìé îàäâéãàýþéêòÿ\nëïšíàñ!"

# ---------------------------------------------------------------------
# 4) Update the saved selection to D19 (matches the commit's new view)
# ---------------------------------------------------------------------
$ws.Range("D19").Select() | Out-Null
